# 2014CoL.xlsx edit: shorten several "City (AltName)" labels to just "City",
# which collapses the "Tehran" shared string into the "Krakow" one (same
# mechanical side effect the original commit introduced), and nudges the
# sheet's saved scroll/selection + default column width.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Cell content edits -----------------------------------------------
# D25  : "Jeddah (Jiddah)"        -> "Jeddah"
# D232 : "Zaragoza (Saragossa)"   -> "Zaragoza"
# D355 : "Krakow (Cracow)"        -> "Krakow"
# D356 : "Tehran"                 -> "Krakow"   (the standalone "Tehran"
#         shared string is dropped, and this cell ends up repointed at the
#         now-renamed "Krakow" string instead)
# D472 : "Tiruchirapalli (Trichy)"-> "Tiruchirapalli"

$ws.Range("D25").Value = "Jeddah"
$ws.Range("D232").Value = "Zaragoza"
$ws.Range("D355").Value = "Krakow"
$ws.Range("D356").Value = "Krakow"
$ws.Range("D472").Value = "Tiruchirapalli"

# --- Sheet view: scroll position + active selection --------------------
$win = $excel.ActiveWindow
$win.ScrollRow = 449
$win.ScrollColumn = 1
$null = $ws.Range("A472").Select()

# --- Default column width ----------------------------------------------
# Original default width attribute 14.4438775510204 -> 14.1734693877551
# (characters). Apply the closest reachable column width across the same
# column span (1 .. 1025) the sheet originally declared.
$ws.Range("A1:AMK1").EntireColumn.ColumnWidth = 13.333333333333334
